$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting the existing rows 16:40 down to 17:41.
$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with the new weekly data point.
$ws.Range("A16").Value = 4
$ws.Range("B16").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C16").Value = 'Los Lagos'
$ws.Range("D16").Value = 44665
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100104
$ws.Range("H16").Value = 'Frutos de pepita'
$ws.Range("I16").Value = 100104001
$ws.Range("J16").Value = 'Granada'
$ws.Range("K16").Value = 'Wonderfull'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 14500
$ws.Range("Q16").Value = '$/caja 14 kilos empedrada'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1036
$ws.Range("T16").Value = 14
